$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings (B1 and C1) to their English equivalents
$ws.Range("B1").Value = "TOTAL_TIME_PER_KLOC_IN_100_COMMITS_WO_FT"
$ws.Range("C1").Value = "TOTAL_TIME_PER_KLOC_IN_100_COMMITS_WITH_FT"

# Update the selected cell in the sheet view to B6
$ws.Range("B6").Select()
